$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.04"
$ws.Range("E2").Value = "'0.19%"
$ws.Range("D3").Value = "'43.91"
$ws.Range("E3").Value = "'-1.54%"
$ws.Range("D4").Value = "'5.513"
$ws.Range("E4").Value = "'-0.64%"
$ws.Range("D5").Value = "'0.08006"
$ws.Range("E5").Value = "'-0.96%"
$ws.Range("D6").Value = "'1.982"
$ws.Range("E6").Value = "'4.27%"
$ws.Range("D7").Value = "'4.298"
$ws.Range("E7").Value = "'-0.83%"
$ws.Range("E8").Value = "'-5.15%"
$ws.Range("D9").Value = "'0.9511"
$ws.Range("E9").Value = "'0.22%"
$ws.Range("D10").Value = "'0.1129"
$ws.Range("E10").Value = "'-4.19%"
$ws.Range("D11").Value = "'0.1862"
$ws.Range("E11").Value = "'-1.76%"
$ws.Range("D12").Value = "'10.50"
$ws.Range("E12").Value = "'25.69%"
$ws.Range("D13").Value = "'0.09883"
$ws.Range("E13").Value = "'-2.34%"
$ws.Range("D14").Value = "'0.04585"
$ws.Range("E14").Value = "'9.68%"
$ws.Range("D15").Value = "'0.1067"
$ws.Range("E15").Value = "'0.20%"
$ws.Range("D16").Value = "'0.001265"
$ws.Range("E16").Value = "'-0.73%"
$ws.Range("D17").Value = "'0.04088"
$ws.Range("E17").Value = "'-3.85%"
$ws.Range("D18").Value = "'0.005861"
$ws.Range("E18").Value = "'-3.49%"
$ws.Range("E19").Value = "'-6.90%"
$ws.Range("D20").Value = "'0.3476"
$ws.Range("E20").Value = "'-0.34%"
$ws.Range("D21").Value = "'0.1408"
$ws.Range("E21").Value = "'2.64%"
$ws.Range("D22").Value = "'0.2548"
$ws.Range("E22").Value = "'-4.32%"
$ws.Range("D23").Value = "'0.001258"
$ws.Range("E23").Value = "'1.88%"
$ws.Range("D24").Value = "'0.004326"
$ws.Range("E24").Value = "'-6.02%"
$ws.Range("D25").Value = "'0.0001159"
$ws.Range("E25").Value = "'-6.03%"
$ws.Range("D26").Value = "'0.0003748"
$ws.Range("E26").Value = "'-6.21%"
$ws.Range("D38").Value = "'0.02550"
$ws.Range("E38").Value = "'-4.24%"
$ws.Range("D39").Value = "'0.05672"
$ws.Range("E39").Value = "'2.08%"
$ws.Range("D40").Value = "'0.007544"
$ws.Range("E40").Value = "'-2.19%"
$ws.Range("D41").Value = "'0.1398"
$ws.Range("E41").Value = "'0.30%"
$ws.Range("D42").Value = "'0.007608"
$ws.Range("E42").Value = "'-32.88%"
$ws.Range("D43").Value = "'0.002014"
$ws.Range("E43").Value = "'-2.19%"
$ws.Range("D44").Value = "'0.008516"
$ws.Range("E44").Value = "'-1.93%"
$ws.Range("D45").Value = "'0.00007115"
$ws.Range("E45").Value = "'-0.10%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'-0.13%"
$ws.Range("E47").Value = "'55.18%"
$ws.Range("D48").Value = "'0.003108"
$ws.Range("E48").Value = "'-9.48%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'-0.13%"
